# Weekly update: insert a new week's "Apio" data (Primera/Segunda) at the
# top of the data block (rows 339-340), shifting the remaining historical
# rows down by two. Everything from row 341 onward is produced automatically
# by the row insert (values/styles move with the rows); we just need to
# populate the two brand-new rows with this week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("339:340").Insert()

# Row 339: Calidad "Primera"
$ws.Range("A339").Value = 6
$ws.Range("B339").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C339").Value = "Metropolitana"
$ws.Range("D339").Value = 44508
$ws.Range("E339").Value = 13
$ws.Range("F339").Value = 100112017
$ws.Range("G339").Value = "Apio"
$ws.Range("H339").Value = "Americana (o)"
$ws.Range("I339").Value = "Primera"
$ws.Range("J339").Value = 1800
$ws.Range("K339").Value = 6000
$ws.Range("L339").Value = 7000
$ws.Range("M339").Value = 6472
$ws.Range("N339").Value = "`$/docena de matas"
$ws.Range("O339").Value = "Región de Coquimbo"
$ws.Range("P339").Value = 1079
$ws.Range("Q339").Value = 6
$ws.Range("R339").Value = "Hortaliza"

# Row 340: Calidad "Segunda"
$ws.Range("A340").Value = 6
$ws.Range("B340").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C340").Value = "Metropolitana"
$ws.Range("D340").Value = 44508
$ws.Range("E340").Value = 13
$ws.Range("F340").Value = 100112017
$ws.Range("G340").Value = "Apio"
$ws.Range("H340").Value = "Americana (o)"
$ws.Range("I340").Value = "Segunda"
$ws.Range("J340").Value = 650
$ws.Range("K340").Value = 5000
$ws.Range("L340").Value = 5000
$ws.Range("M340").Value = 5000
$ws.Range("N340").Value = "`$/docena de matas"
$ws.Range("O340").Value = "Región de Coquimbo"
$ws.Range("P340").Value = 833
$ws.Range("Q340").Value = 6
$ws.Range("R340").Value = "Hortaliza"
